# Auto-generated edit script: apply numeric corrections to Ixion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 277.4
$ws.Range("I28").Value = 208.13333
$ws.Range("J28").Value = 485.2
$ws.Range("K28").Value = 208.13333
$ws.Range("L28").Value = 485.2
$ws.Range("M28").Value = 276.86667
$ws.Range("N28").Value = -1455.2

$ws.Range("H116").Value = 8313.157999999999
$ws.Range("I116").Value = 8608.333000000001
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 8608.333000000001
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -5166.333000000001
$ws.Range("N116").Value = -9884

$ws.Range("H129").Value = 1383.7646
$ws.Range("I129").Value = 832
$ws.Range("J129").Value = 1613.6666
$ws.Range("K129").Value = 2496
$ws.Range("L129").Value = 4840.9998
$ws.Range("M129").Value = 2504
$ws.Range("N129").Value = -14840.9998

$ws.Range("H137").Value = 1539.2128
$ws.Range("I137").Value = 1200.9524
$ws.Range("J137").Value = 4380.6
$ws.Range("K137").Value = 3602.857199999999
$ws.Range("L137").Value = 13141.8
$ws.Range("M137").Value = -1052.857199999999
$ws.Range("N137").Value = -18241.8

$ws.Range("H138").Value = 2691.879
$ws.Range("I138").Value = 1506.9032
$ws.Range("J138").Value = 3741.4285
$ws.Range("K138").Value = 4520.7096
$ws.Range("L138").Value = 11224.2855
$ws.Range("M138").Value = 619.2903999999999
$ws.Range("N138").Value = -21504.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7540.69
$ws.Range("I32").Value = 5299.6465
$ws.Range("J32").Value = 17749.889
$ws.Range("K32").Value = 5299.6465
$ws.Range("L32").Value = 17749.889
$ws.Range("M32").Value = -5012.6465
$ws.Range("N32").Value = -18323.889

$ws.Range("H45").Value = 111267.48
$ws.Range("I45").Value = 150583.42
$ws.Range("J45").Value = 1182.8
$ws.Range("K45").Value = 150583.42
$ws.Range("L45").Value = 1182.8
$ws.Range("M45").Value = -150206.42
$ws.Range("N45").Value = -1936.8

$ws.Range("H97").Value = 628.6667
$ws.Range("I97").Value = 374.6154
$ws.Range("J97").Value = 1041.5
$ws.Range("K97").Value = 374.6154
$ws.Range("L97").Value = 1041.5
$ws.Range("M97").Value = 121.3846
$ws.Range("N97").Value = -2033.5

$ws.Range("H122").Value = 1168109.5
$ws.Range("I122").Value = 1223676.6
$ws.Range("K122").Value = 3671029.8
$ws.Range("M122").Value = -3668579.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 320.2857
$ws.Range("J64").Value = 335.9091
$ws.Range("L64").Value = 335.9091
$ws.Range("N64").Value = -785.9091000000001

$ws.Range("H67").Value = 320.2857
$ws.Range("J67").Value = 335.9091
$ws.Range("L67").Value = 335.9091
$ws.Range("N67").Value = -1895.9091

$ws.Range("H105").Value = 1783.625
$ws.Range("I105").Value = 1693.8
$ws.Range("K105").Value = 1693.8
$ws.Range("M105").Value = 53.20000000000005

$ws.Range("H137").Value = 51720
$ws.Range("J137").Value = 44435
$ws.Range("L137").Value = 44435
$ws.Range("N137").Value = -54635

$ws.Range("H138").Value = 43812
$ws.Range("J138").Value = 43812
$ws.Range("L138").Value = 43812
$ws.Range("N138").Value = -54092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9917841
$ws.Range("I31").Value = 1899.9286
$ws.Range("J31").Value = 13571082
$ws.Range("K31").Value = 1899.9286
$ws.Range("L31").Value = 13571082
$ws.Range("N31").Value = -13571672
$ws.Range("M31").Value = -1604.9286

$ws.Range("H34").Value = 9917841
$ws.Range("I34").Value = 1899.9286
$ws.Range("J34").Value = 13571082
$ws.Range("K34").Value = 1899.9286
$ws.Range("L34").Value = 13571082
$ws.Range("N34").Value = -13571486
$ws.Range("M34").Value = -1697.9286

$ws.Range("H58").Value = 6311335.5
$ws.Range("I58").Value = 8335068.5
$ws.Range("J58").Value = 1252003.6
$ws.Range("K58").Value = 8335068.5
$ws.Range("L58").Value = 1252003.6
$ws.Range("M58").Value = -8334865.5
$ws.Range("N58").Value = -1252409.6

$ws.Range("H132").Value = 10529473
$ws.Range("I132").Value = 18183916
$ws.Range("J132").Value = 4613.5
$ws.Range("K132").Value = 54551748
$ws.Range("L132").Value = 13840.5
$ws.Range("M132").Value = -54549218
$ws.Range("N132").Value = -18900.5

$ws.Range("H134").Value = 9540812
$ws.Range("I134").Value = 16670913
$ws.Range("J134").Value = 628186.3
$ws.Range("K134").Value = 50012739
$ws.Range("L134").Value = 1884558.9
$ws.Range("M134").Value = -50010204
$ws.Range("N134").Value = -1889628.9

$ws.Range("H136").Value = 6311335.5
$ws.Range("I136").Value = 8335068.5
$ws.Range("J136").Value = 1252003.6
$ws.Range("K136").Value = 25005205.5
$ws.Range("L136").Value = 3756010.8
$ws.Range("M136").Value = -25002655.5
$ws.Range("N136").Value = -3761110.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1696069
$ws.Range("J131").Value = 1335.7872
$ws.Range("L131").Value = 4007.3616
$ws.Range("N131").Value = -14087.3616

$ws.Range("H137").Value = 19396.334
$ws.Range("I137").Value = 9286
$ws.Range("J137").Value = 44672.168
$ws.Range("K137").Value = 27858
$ws.Range("L137").Value = 134016.504
$ws.Range("M137").Value = -22758
$ws.Range("N137").Value = -144216.504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4555.7915
$ws.Range("I70").Value = 4264.5
$ws.Range("J70").Value = 5138.375
$ws.Range("K70").Value = 4264.5
$ws.Range("L70").Value = 5138.375
$ws.Range("M70").Value = -3994.5
$ws.Range("N70").Value = -5678.375

$ws.Range("H73").Value = 4555.7915
$ws.Range("I73").Value = 4264.5
$ws.Range("J73").Value = 5138.375
$ws.Range("K73").Value = 4264.5
$ws.Range("L73").Value = 5138.375
$ws.Range("M73").Value = -3328.5
$ws.Range("N73").Value = -7010.375

$ws.Range("H102").Value = 1965.8096
$ws.Range("I102").Value = 1576
$ws.Range("J102").Value = 2940.3333
$ws.Range("K102").Value = 1576
$ws.Range("L102").Value = 2940.3333
$ws.Range("M102").Value = 46
$ws.Range("N102").Value = -6184.3333

$ws.Range("H113").Value = 41667670
$ws.Range("I113").Value = 71429360
$ws.Range("J113").Value = 1309
$ws.Range("K113").Value = 71429360
$ws.Range("L113").Value = 1309
$ws.Range("M113").Value = -71427190
$ws.Range("N113").Value = -5649

$ws.Range("H122").Value = 48745190
$ws.Range("I122").Value = 66552290
$ws.Range("J122").Value = 25002400
$ws.Range("K122").Value = 199656870
$ws.Range("L122").Value = 75007200
$ws.Range("M122").Value = -199654420
$ws.Range("N122").Value = -75012100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1173
$ws.Range("I113").Value = 1257.125
$ws.Range("K113").Value = 3771.375
$ws.Range("M113").Value = -1601.375

$ws.Range("H122").Value = 1409
$ws.Range("I122").Value = 1390.8
$ws.Range("K122").Value = 4172.4
$ws.Range("M122").Value = -1722.4

$ws.Range("H132").Value = 1734.7354
$ws.Range("I132").Value = 749.3333
$ws.Range("J132").Value = 2843.3125
$ws.Range("K132").Value = 2247.9999
$ws.Range("L132").Value = 8529.9375
$ws.Range("M132").Value = 282.0001000000002
$ws.Range("N132").Value = -13589.9375
